# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" sheet (cloned from "2022-Q1" so it keeps the
#    same layout/styles) positioned right before "2022-Q1".
# 2. Fill it in with the new quarter's fund data.
# 3. Insert a matching new row into the "总计" (totals) summary sheet and
#    renumber the trailing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2. New "2022-Q3" worksheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

# D2:G3 hold text-formatted numbers (percentages / amounts stored as
# strings, same as the rest of the workbook) - force text so Excel does
# not silently re-interpret them as numeric values.
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("D2").Value = "10.53"
$q3.Range("E2").Value = "96.33"
$q3.Range("F2").Value = "9.12"
$q3.Range("G2").Value = "0.9603"

$q3.Range("D3").Value = "0.02"
$q3.Range("E3").Value = "42.74"
$q3.Range("F3").Value = "3.00"
$q3.Range("G3").Value = "0.0006"
$q3.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 3. "总计" summary sheet gets a new row for 2022-Q3
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.96

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
